$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update amounts
$ws.Range("C8").Value = 6.44
$ws.Range("C9").Value = 16.36

# Capture the original hyperlink target URLs in their original (E2..E10)
# order. Hyperlink.Address isn't populated for links loaded from disk in
# this engine, but the URL happens to equal the cell's own displayed text
# here, so pull it from there while the links are still on their original
# cells (this does NOT change any cell content).
$originalTargets = @()
for ($r = 2; $r -le 10; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $originalTargets += , $cell.Text
}

$ws.Hyperlinks.Delete()

# Re-attach the same nine target URLs to the cells, but in the new
# (mostly-reversed) arrangement, and in this element order, matching the
# edited workbook: the hyperlink that used to sit on E2 (1st target) now
# sits on E10, the one that used to sit on E3 (2nd target) now sits on E9,
# etc. Cell contents (the visible text in column E) are left untouched.
$newOrder = @(
    @{ Ref = "E10"; TargetIndex = 0 },
    @{ Ref = "E9";  TargetIndex = 1 },
    @{ Ref = "E8";  TargetIndex = 2 },
    @{ Ref = "E7";  TargetIndex = 3 },
    @{ Ref = "E6";  TargetIndex = 4 },
    @{ Ref = "E5";  TargetIndex = 5 },
    @{ Ref = "E3";  TargetIndex = 6 },
    @{ Ref = "E4";  TargetIndex = 7 },
    @{ Ref = "E2";  TargetIndex = 8 }
)
foreach ($entry in $newOrder) {
    $ws.Hyperlinks.Add($ws.Range($entry.Ref), $originalTargets[$entry.TargetIndex]) | Out-Null
}

# Hyperlinks.Add() re-derives a (slightly different, but equivalent) cell
# style record for the hyperlinked cells; restore the original named
# "Hyperlink" cell style so the cells keep referencing the same style as
# before instead of a newly duplicated one.
$ws.Range("E2:E10").Style = "Hyperlink"

# Move the active selection to C9 (single cell)
$ws.Range("C9").Select() | Out-Null
